$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 4) mirroring the content pattern of the
# existing rows (row 2 / row 3).
$ws.Range("A4").Value = "FUL_Transmittals_ActionOverDue_New_ChangeNote"
$ws.Range("B4").Value = "Creates a new Transmittal of  Change Note and validate the count in Action Overdue"
$ws.Range("C4").Value = "N"
$ws.Range("D4").Value = "Y"
$ws.Range("F4").Value = "Sprint1"

# Match formatting of the previous data row (copy formats only, values
# already set above)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B3:F3").Copy()
$ws.Range("B4:F4").PasteSpecial(-4122)

# Extend the data validation lists down to the new row (recreate, in the
# same order, so the new row is covered by the sqref ranges)
$ws.Range("F2:F3").Validation.Delete()
$ws.Range("C2:D3").Validation.Delete()
$ws.Range("C2:D4").Validation.Add(3, 1, 1, '"Y,N"')
$ws.Range("F2:F4").Validation.Add(3, 1, 1, '"Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10"')

# Update the active selection to mirror the diff
$ws.Range("D2:D4").Select()
